$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Defaults")

# Remember the existing mailto hyperlinks (cell -> target) before the row
# shift below invalidates their anchoring; the engine does not re-anchor
# hyperlink ranges automatically when rows are deleted.
$linkTargets = @(
    "mailto:Manchesteret@justice.gov.uk",
    "mailto:glasgowet@justice.gov.uk",
    "mailto:aberdeenet@justice.gov.uk",
    "mailto:dundeeet@justice.gov.uk",
    "mailto:edinburghet@justice.gov.uk",
    "mailto:bristolet@justice.gov.uk",
    "mailto:LeedsET@justice.gov.uk",
    "mailto:londoncentralet@hmcts.gsi.gov.uk",
    "mailto:eastlondon@justice.gov.uk",
    "mailto:londonsouthet@hmcts.gsi.gov.uk",
    "mailto:e.midlandseastet@justice.gov.uk",
    "mailto:MidlandsWestET@justice.gov.uk",
    "mailto:newcastleet@hmcts.gsi.gov.uk",
    "mailto:cardiffet@justice.gov.uk",
    "mailto:watfordet@justice.gov.uk"
)
$linkRows = @(11, 20, 29, 38, 45, 53, 61, 70, 78, 87, 94, 102, 111, 118, 127)

$ws.Hyperlinks.Delete()

# Combine the two Glasgow address lines (AddressLine2 + AddressLine3) into a
# single AddressLine2 value, then remove the now-redundant AddressLine3 row
# (row 14: "tribunalGlasgowAddressLine3" / "20 York Street").
$ws.Range("B13").Value = "3 Atlantic Quay, 20 York Street"
$ws.Rows.Item(14).Delete()

# Re-create the hyperlinks at their shifted rows (every row from 14 down
# moved up by one). Restore the plain "Normal" cell style afterwards since
# the original cells were not using the auto-applied "Hyperlink" style.
for ($i = 0; $i -lt $linkRows.Count; $i++) {
    $row = $linkRows[$i]
    if ($row -ge 14) {
        $row = $row - 1
    }
    $target = $ws.Range("B$row")
    $ws.Hyperlinks.Add($target, $linkTargets[$i], "", "", "")
    $target.Style = "Normal"
}

# Update the active selection to match the new layout.
$ws.Range("B18").Select()
